# Regenerate the "K" column (column G) values on the active worksheet.
# The source data feeding this sheet was regenerated upstream (K instead of
# Strike#), so here we just write out the newly computed values for each
# existing data row (rows 2-66), leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 4
    11 = 3
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 2
    20 = 2
    21 = 2
    22 = 0
    23 = 0
    25 = 2
    26 = 4
    27 = 2
    28 = 4
    29 = 2
    30 = 1
    31 = 3
    32 = 0
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 0
    49 = 0
    50 = 3
    51 = 0
    52 = 1
    53 = 1
    54 = 3
    55 = 0
    56 = 3
    57 = 2
    58 = 3
    59 = 1
    60 = 1
    61 = 2
    62 = 1
    63 = 2
    64 = 1
    65 = 2
    66 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
